# Append two new data rows (176, 177) to Sheet1, matching the source data
# feed pattern: one row for each station ("四方坪站充电量(kw)" / "高岭站充电量(kw)")
# for date-serial 45988 (2025-11-27), followed by 24 hourly values (columns C:Z).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row176 = @(676.25299999999993, 836.35399999999993, 290.11199999999997, 345.55099999999999, 390.39300000000003, 661.13499999999999, 407.92399999999992, 187.624, 203.24299999999999, 93.02000000000001, 92.51, 235.21299999999999, 971.72700000000009, 1064.3489999999999, 280.83399999999995, 294.75799999999998, 438.09000000000003, 82.544000000000011, 91.58, 177.91400000000002, 72.161000000000001, 174.45999999999998, 91.899999999999991, 65.16)

$row177 = @(504.50799999999987, 381.69200000000001, 184.90699999999998, 72.238, 169.696, 364.05500000000001, 43.14, 67.012, 297.10700000000003, 240.47800000000001, 179.923, 359.221, 351.75299999999999, 458.20500000000004, 299.08799999999997, 269.59600000000006, 129.59800000000001, 76.906000000000006, 3.68, 115.506, 0, 32.671999999999997, 36.280999999999999, 0.36)

# Row 176: 四方坪站充电量(kw)
$ws.Cells.Item(176, 1).Value = 45988
$ws.Cells.Item(176, 2).Value = "四方坪站充电量(kw)"
for ($i = 0; $i -lt $row176.Length; $i++) {
    $ws.Cells.Item(176, 3 + $i).Value = $row176[$i]
}

# Row 177: 高岭站充电量(kw)
$ws.Cells.Item(177, 1).Value = 45988
$ws.Cells.Item(177, 2).Value = "高岭站充电量(kw)"
for ($i = 0; $i -lt $row177.Length; $i++) {
    $ws.Cells.Item(177, 3 + $i).Value = $row177[$i]
}

# Match the style of the rows directly above (date format on col A, number
# format on C:Z) by copying formats down from row 175.
$ws.Range("A175:Z175").Copy() | Out-Null
$ws.Range("A176:Z177").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Final selection moves to H189, matching the new cursor position after entry.
$ws.Range("H189").Select() | Out-Null
